{"js": "// Remove the thin gray \"\u2500\u2500\u2500\u2500\u2500\u2500\" separator paragraphs and the tiny empty\n// spacer paragraphs (pPr spacing w:before=\"40\", no text) that precede\n// them after a code-sample table, while leaving every other paragraph \u2014\n// including all images/drawings \u2014 untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Need paragraph-level spacing (in points: 40 twips = 2pt, 120 twips = 6pt)\n// plus whether each paragraph lives inside a table cell, so we never touch\n// table content.\nfor (const p of paragraphs.items) {\n  p.paragraphFormat.load(\"spaceBefore,spaceAfter\");\n}\nawait context.sync();\n\nconst nullObjChecks = paragraphs.items.map((p) => p.parentTableCellOrNullObject);\nfor (const tc of nullObjChecks) {\n  tc.load(\"isNullObject\");\n}\nawait context.sync();\n\nconst separatorPattern = /^[\\u2500]+$/; // the \"\u2500\" box-drawing character repeated\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const inTable = !nullObjChecks[i].isNullObject;\n  if (inTable) continue; // never touch table contents (preserves the code listings)\n\n  const text = p.text;\n  const isSeparatorLine = text.length > 0 && separatorPattern.test(text);\n  const isEmptySpacer =\n    text === \"\" &&\n    p.paragraphFormat.spaceBefore === 2 &&\n    p.paragraphFormat.spaceAfter !== 6;\n\n  if (isSeparatorLine || isEmptySpacer) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the thin gray \"\u2500\u2500\u2500\u2500\u2500\u2500\" separator paragraphs and the tiny empty\n# spacer paragraphs (spacing before=2pt / 40 twips, no text) that sit\n# right before them after a code-sample table, while leaving every other\n# paragraph -- including all images -- untouched.\n#\n# Walk the paragraph collection from the END towards the start and delete\n# matches by re-fetching Item($i) fresh each time. Doing this back-to-front\n# means every paragraph we still need to inspect keeps its original index,\n# so earlier deletions never disturb later lookups.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n    # Strip the trailing paragraph-mark / cell-mark control characters\n    # (Chr(13) / Chr(7)) before comparing the visible text.\n    $trimmed = $text -replace \"[\\r\\a]+$\", \"\"\n\n    $spaceBefore = $p.Format.SpaceBefore\n    $spaceAfter = $p.Format.SpaceAfter\n\n    $isSeparatorLine = $trimmed -match \"^[\\u2500]+$\"\n    $isEmptySpacer = ($trimmed -eq \"\") -and ($spaceBefore -eq 2) -and ($spaceAfter -ne 6)\n\n    if ($isSeparatorLine -or $isEmptySpacer) {\n        $p.Range.Delete()\n    }\n}\n"}
